$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the "Femacal de La Calera - Alcachofa" price series.
# Two new weekly observations are inserted at the top of the data block
# (rows 169-170), which pushes every previously-recorded observation down
# by two rows, and the two oldest observations that fall off the bottom
# of the historical window are appended as brand-new rows 182-183.
$rows = @(
    [PSCustomObject]@{ Row=169; D=44461; H='Argentina(o)'; I='Primera'; J=65; K=10000; L=10000; M=10000; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=200; Q=50 }
    [PSCustomObject]@{ Row=170; D=44461; H='Española'; I='Extra'; J=95; K=11500; L=12000; M=11737; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=391; Q=30 }
    [PSCustomObject]@{ Row=171; D=44357; H='Española'; I='Extra'; J=115; K=14000; L=15000; M=14478; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=483; Q=30 }
    [PSCustomObject]@{ Row=172; D=44162; H='Española'; I='Primera'; J=3800; K=250; L=250; M=250; N='$/unidad'; O='Llay Llay'; P=250; Q=1 }
    [PSCustomObject]@{ Row=173; D=44410; H='Argentina(o)'; I='Primera'; J=250; K=13000; L=13500; M=13240; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=265; Q=50 }
    [PSCustomObject]@{ Row=174; D=44410; H='Española'; I='Extra'; J=120; K=13500; L=14000; M=13750; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=458; Q=30 }
    [PSCustomObject]@{ Row=175; D=44411; H='Argentina(o)'; I='Primera'; J=110; K=12500; L=13000; M=12727; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=255; Q=50 }
    [PSCustomObject]@{ Row=176; D=44411; H='Española'; I='Extra'; J=60; K=13500; L=13500; M=13500; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=450; Q=30 }
    [PSCustomObject]@{ Row=177; D=44411; H='Española'; I='Primera'; J=60; K=13000; L=13000; M=13000; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=433; Q=30 }
    [PSCustomObject]@{ Row=178; D=44376; H='Argentina(o)'; I='Primera'; J=60; K=15000; L=15000; M=15000; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=300; Q=50 }
    [PSCustomObject]@{ Row=179; D=44376; H='Española'; I='Extra'; J=105; K=16000; L=16500; M=16238; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=541; Q=30 }
    [PSCustomObject]@{ Row=180; D=44358; H='Argentina(o)'; I='Primera'; J=110; K=14000; L=14500; M=14227; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=285; Q=50 }
    [PSCustomObject]@{ Row=181; D=44358; H='Española'; I='Extra'; J=60; K=14000; L=14000; M=14000; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=467; Q=30 }
    [PSCustomObject]@{ Row=182; D=44425; H='Argentina(o)'; I='Primera'; J=110; K=12500; L=13000; M=12773; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=255; Q=50 }
    [PSCustomObject]@{ Row=183; D=44425; H='Española'; I='Extra'; J=125; K=12000; L=12500; M=12240; N='$/caja 30 unidades'; O='Provincia de Limarí'; P=408; Q=30 }
)

# Columns that never change for this sub-sheet (every row describes the
# same market / region / category / classification combination).
$constCols = @{
    A = 3
    B = 'Femacal de La Calera'
    C = 'Coquimbo'
    E = 5
    F = 100112013
    G = 'Alcachofa'
    R = 'Hortaliza'
}

$dateFormat = $ws.Range("D169").NumberFormat

foreach ($r in $rows) {
    $row = $r.Row

    foreach ($col in $constCols.Keys) {
        $ws.Range("$col$row").Value = $constCols[$col]
    }

    $ws.Range("D$row").Value = $r.D
    $ws.Range("D$row").NumberFormat = $dateFormat
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
}
